$d = $word.ActiveDocument

# 1. Contact info line: drop phone number + street address,
#    add "Yardley Bucks County Area" instead.
$r = $d.Content
$r.Find.Execute(
    "267-469-1210 | 3kw109@gmail.com | 548 Stevens Rd, Morrisville, PA 19067",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "| 3kw109@gmail.com | ", 2)
$r.Collapse(0)
$r.InsertAfter("Yardley Bucks County Area")

# 2. "Majoring in Statistics" -> "Applied Statistics"
$r = $d.Content
$r.Find.Execute(
    "Majoring in", $true, $false, $false, $false, $false, $true, 1, $false,
    "Applied", 2)

# 3. "Mathematical Statistics 2," -> "Mathematical Statistics,"
$r = $d.Content
$r.Find.Execute(
    "Mathematical Statistics 2,", $true, $false, $false, $false, $false, $true, 1, $false,
    "Mathematical Statistics,", 2)

# 4. "Expected graduation: " -> "Graduation: "
$r = $d.Content
$r.Find.Execute(
    "Expected graduation: ", $true, $false, $false, $false, $false, $true, 1, $false,
    "Graduation: ", 2)
